$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new rows (8-11) of keyword data for a new title, matching
# the existing table layout (ASIN, Title, Keyword).
$asin = "B01531YHYU"
$title = "Project Return Fire"

$ws.Range("A8").Value = $asin
$ws.Range("B8").Value = $title
$ws.Range("C8").Value = "Action Adventure"

$ws.Range("A9").Value = $asin
$ws.Range("B9").Value = $title
$ws.Range("C9").Value = "Time Travel Action Adventure"

$ws.Range("A10").Value = $asin
$ws.Range("B10").Value = $title
$ws.Range("C10").Value = "Science Fiction Action Adventure"

$ws.Range("A11").Value = $asin
$ws.Range("B11").Value = $title
$ws.Range("C11").Value = "Science Fiction"

# Match the author's final selection/navigation position after entering data.
$ws.Range("B16").Select()
